$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.509.40"
$ws.Cells.Item(2, 5).Value = "  +4.35%  "
$ws.Cells.Item(3, 4).Value = "1.841.95"
$ws.Cells.Item(3, 5).Value = "  +3.73%  "
$ws.Cells.Item(4, 4).Value = "'1.030"
$ws.Cells.Item(4, 5).Value = "  +2.87%  "
$ws.Cells.Item(5, 4).Value = "'319.21"
$ws.Cells.Item(5, 5).Value = "  +4.70%  "
$ws.Cells.Item(6, 4).Value = "'1.028"
$ws.Cells.Item(7, 5).Value = "  +3.23%  "
$ws.Cells.Item(8, 4).Value = "'0.3733"
$ws.Cells.Item(8, 5).Value = "  +3.76%  "
$ws.Cells.Item(9, 4).Value = "'0.07390"
$ws.Cells.Item(9, 5).Value = "  +3.22%  "
$ws.Cells.Item(10, 4).Value = "'0.8751"
$ws.Cells.Item(10, 5).Value = "  +4.64%  "
$ws.Cells.Item(11, 4).Value = "'21.41"
$ws.Cells.Item(11, 5).Value = "  +4.76%  "
$ws.Cells.Item(12, 4).Value = "1.858.42"
$ws.Cells.Item(12, 5).Value = "  +4.85%  "
$ws.Cells.Item(13, 5).Value = "  +4.54%  "
$ws.Cells.Item(14, 4).Value = "'6.691"
$ws.Cells.Item(14, 5).Value = "  +3.93%  "
$ws.Cells.Item(15, 4).Value = "'0.07161"
$ws.Cells.Item(15, 5).Value = "  +3.76%  "
$ws.Cells.Item(16, 4).Value = "'82.60"
$ws.Cells.Item(16, 5).Value = "  +4.47%  "
$ws.Cells.Item(17, 4).Value = "'1.032"
$ws.Cells.Item(17, 5).Value = "  +3.10%  "
$ws.Cells.Item(18, 4).Value = "'0.000009025"
$ws.Cells.Item(18, 5).Value = "  +4.35%  "
$ws.Cells.Item(19, 5).Value = "  +2.60%  "
$ws.Cells.Item(20, 4).Value = "'15.40"
$ws.Cells.Item(20, 5).Value = "  +3.39%  "
$ws.Cells.Item(21, 4).Value = "27.527.76"
$ws.Cells.Item(21, 5).Value = "  +4.39%  "
$ws.Cells.Item(22, 4).Value = "'5.229"
$ws.Cells.Item(22, 5).Value = "  +3.05%  "
$ws.Cells.Item(23, 4).Value = "'11.29"
$ws.Cells.Item(23, 5).Value = "  +3.68%  "
$ws.Cells.Item(24, 4).Value = "2.066.96"
$ws.Cells.Item(24, 5).Value = "  +4.11%  "
$ws.Cells.Item(25, 4).Value = "'156.99"
$ws.Cells.Item(25, 5).Value = "  +3.66%  "
$ws.Cells.Item(26, 4).Value = "'1.912"
$ws.Cells.Item(26, 5).Value = "  +6.18%  "
$ws.Cells.Item(27, 4).Value = "'18.64"
$ws.Cells.Item(27, 5).Value = "  +3.81%  "
$ws.Cells.Item(28, 4).Value = "'5.263"
$ws.Cells.Item(28, 5).Value = "  +3.64%  "
$ws.Cells.Item(29, 5).Value = "  +5.40%  "
$ws.Cells.Item(30, 4).Value = "'116.20"
$ws.Cells.Item(30, 5).Value = "  +1.60%  "
$ws.Cells.Item(31, 5).Value = "  +2.93%  "
$ws.Cells.Item(32, 5).Value = "  +7.60%  "
$ws.Cells.Item(33, 4).Value = "'0.7624"
$ws.Cells.Item(33, 5).Value = "  +4.84%  "
$ws.Cells.Item(34, 5).Value = "  +3.87%  "
$ws.Cells.Item(35, 5).Value = "  +5.25%  "
$ws.Cells.Item(37, 4).Value = "'1.148"
$ws.Cells.Item(37, 5).Value = "  +5.58%  "
$ws.Cells.Item(38, 4).Value = "'0.01968"
$ws.Cells.Item(38, 5).Value = "  +4.64%  "
$ws.Cells.Item(39, 4).Value = "'0.05247"
$ws.Cells.Item(39, 5).Value = "  +2.71%  "
$ws.Cells.Item(40, 4).Value = "'0.5173"
$ws.Cells.Item(40, 5).Value = "  +5.15%  "
$ws.Cells.Item(41, 4).Value = "'2.777"
$ws.Cells.Item(41, 5).Value = "  +7.01%  "
$ws.Cells.Item(42, 4).Value = "'0.1663"
$ws.Cells.Item(42, 5).Value = "  +3.57%  "
$ws.Cells.Item(43, 4).Value = "'6.578"
$ws.Cells.Item(43, 5).Value = "  +4.09%  "
$ws.Cells.Item(44, 4).Value = "'8.494"
$ws.Cells.Item(44, 5).Value = "  +5.90%  "
$ws.Cells.Item(45, 2).Value = "Quant"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(45, 4).Value = "'108.98"
$ws.Cells.Item(45, 5).Value = "  +4.48%  "
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).Value = "'10.63"
$ws.Cells.Item(46, 5).Value = "  +4.71%  "
$ws.Cells.Item(47, 2).Value = "PaxDollar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(47, 4).Value = "'1.031"
$ws.Cells.Item(47, 5).Value = "  +3.01%  "
$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 4).Value = "'1.702"
$ws.Cells.Item(48, 5).Value = "  +4.40%  "
$ws.Cells.Item(49, 2).Value = "Decentraland"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(49, 4).Value = "'0.4637"
$ws.Cells.Item(49, 5).Value = "  +4.35%  "
$ws.Cells.Item(50, 2).Value = "RenderToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(50, 4).Value = "'1.894"
$ws.Cells.Item(50, 5).Value = "  +10.10%  "
$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.06333"
$ws.Cells.Item(51, 5).Value = "  +2.70%  "
